$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4
$ws.Range("C2").Value = "flower/flower086.png"
$ws.Range("D2").Value = "reisen"
$ws.Range("E2").Value = "flower"

$ws.Range("B3").Value = 86
$ws.Range("C3").Value = "flower/flower101.png"
$ws.Range("D3").Value = "deuten"
$ws.Range("E3").Value = "flower"

$ws.Range("B4").Value = 101
$ws.Range("C4").Value = "face/face079.png"
$ws.Range("D4").Value = "rufen"
$ws.Range("E4").Value = "face"

$ws.Range("B5").Value = 29
$ws.Range("C5").Value = "flower/flower070.png"
$ws.Range("D5").Value = "öffnen"
$ws.Range("E5").Value = "flower"

$ws.Range("B6").Value = 16
$ws.Range("C6").Value = "face/face085.png"
$ws.Range("D6").Value = "kranken"
$ws.Range("E6").Value = "face"

$ws.Range("B7").Value = 118
$ws.Range("C7").Value = "face/face100.png"
$ws.Range("D7").Value = "parken"
$ws.Range("E7").Value = "face"

$ws.Range("B8").Value = 72
$ws.Range("C8").Value = "flower/flower113.png"
$ws.Range("D8").Value = "planen"
$ws.Range("E8").Value = "flower"

$ws.Range("B9").Value = 63
$ws.Range("C9").Value = "flower/flower069.png"
$ws.Range("D9").Value = "holen"
$ws.Range("E9").Value = "flower"

$ws.Range("B10").Value = 47
$ws.Range("C10").Value = "flower/flower083.png"
$ws.Range("D10").Value = "trotzen"
$ws.Range("E10").Value = "flower"

$ws.Range("B11").Value = 94
$ws.Range("C11").Value = "face/face112.png"
$ws.Range("D11").Value = "wehen"
$ws.Range("E11").Value = "face"

$ws.Range("B12").Value = 105
$ws.Range("C12").Value = "face/face096.png"
$ws.Range("D12").Value = "frischen"
$ws.Range("E12").Value = "face"

$ws.Range("B13").Value = 84
$ws.Range("C13").Value = "face/face076.png"
$ws.Range("D13").Value = "tollen"
$ws.Range("E13").Value = "face"

$ws.Range("B14").Value = 11
$ws.Range("C14").Value = "face/face073.png"
$ws.Range("D14").Value = "prüfen"
$ws.Range("E14").Value = "face"

$ws.Range("B15").Value = 123
$ws.Range("C15").Value = "face/face098.png"
$ws.Range("D15").Value = "nullen"
$ws.Range("E15").Value = "face"

$ws.Range("B16").Value = 15
$ws.Range("C16").Value = "face/face078.png"
$ws.Range("D16").Value = "piepen"
$ws.Range("E16").Value = "face"

$ws.Range("B17").Value = 51
$ws.Range("C17").Value = "face/face091.png"
$ws.Range("D17").Value = "kennen"
$ws.Range("E17").Value = "face"

$ws.Range("B18").Value = 7
$ws.Range("C18").Value = "flower/flower065.png"
$ws.Range("D18").Value = "bauen"
$ws.Range("E18").Value = "flower"

$ws.Range("B19").Value = 108
$ws.Range("C19").Value = "face/face103.png"
$ws.Range("D19").Value = "meinen"
$ws.Range("E19").Value = "face"

$ws.Range("B20").Value = 109
$ws.Range("C20").Value = "flower/flower073.png"
$ws.Range("D20").Value = "narren"
$ws.Range("E20").Value = "flower"

$ws.Range("B21").Value = 3
$ws.Range("C21").Value = "face/face083.png"
$ws.Range("D21").Value = "achten"
$ws.Range("E21").Value = "face"

$ws.Range("B22").Value = 73
$ws.Range("C22").Value = "flower/flower103.png"
$ws.Range("D22").Value = "rechnen"
$ws.Range("E22").Value = "flower"

$ws.Range("B23").Value = 40
$ws.Range("C23").Value = "face/face093.png"
$ws.Range("D23").Value = "mögen"
$ws.Range("E23").Value = "face"

$ws.Range("B24").Value = 76
$ws.Range("C24").Value = "face/face107.png"
$ws.Range("D24").Value = "sparen"
$ws.Range("E24").Value = "face"

$ws.Range("B25").Value = 60
$ws.Range("C25").Value = "flower/flower078.png"
$ws.Range("D25").Value = "heben"
$ws.Range("E25").Value = "flower"

$ws.Range("B26").Value = 28
$ws.Range("C26").Value = "flower/flower067.png"
$ws.Range("D26").Value = "lassen"
$ws.Range("E26").Value = "flower"

$ws.Range("B27").Value = 69
$ws.Range("C27").Value = "flower/flower099.png"
$ws.Range("D27").Value = "stoppen"
$ws.Range("E27").Value = "flower"

$ws.Range("B28").Value = 14
$ws.Range("C28").Value = "flower/flower080.png"
$ws.Range("D28").Value = "lügen"
$ws.Range("E28").Value = "flower"

$ws.Range("B29").Value = 88
$ws.Range("C29").Value = "flower/flower094.png"
$ws.Range("D29").Value = "heißen"
$ws.Range("E29").Value = "flower"

$ws.Range("B30").Value = 43
$ws.Range("C30").Value = "face/face066.png"
$ws.Range("D30").Value = "passen"
$ws.Range("E30").Value = "face"

$ws.Range("B31").Value = 111
$ws.Range("C31").Value = "flower/flower114.png"
$ws.Range("D31").Value = "dienen"
$ws.Range("E31").Value = "flower"

$ws.Range("B32").Value = 61
$ws.Range("C32").Value = "flower/flower066.png"
$ws.Range("D32").Value = "binden"
$ws.Range("E32").Value = "flower"

$ws.Range("B33").Value = 59
$ws.Range("C33").Value = "face/face102.png"
$ws.Range("D33").Value = "hassen"
$ws.Range("E33").Value = "face"
